$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 8.600386905322791
$ws.Range("C2").Value = 4.961210941492835
$ws.Range("D2").Value = 4.799110630911748
$ws.Range("F2").Value = 22.00976878969509
$ws.Range("G2").Value = 3.621657502333183
$ws.Range("I2").Value = 18.96588030718569
$ws.Range("K2").Value = 8.350989176995258
$ws.Range("N2").Value = 17.48797667798502
$ws.Range("O2").Value = 19.85960241583157
$ws.Range("B3").Value = 8.280044250274301
$ws.Range("C3").Value = 4.77783824328683
$ws.Range("D3").Value = 4.706033055650416
$ws.Range("F3").Value = 22.06867399777129
$ws.Range("G3").Value = 3.623041583841587
$ws.Range("I3").Value = 19.05714943363522
$ws.Range("K3").Value = 8.133181791922043
$ws.Range("N3").Value = 17.53437867516515
$ws.Range("O3").Value = 19.94265453520549
$ws.Range("B4").Value = 8.077764362241856
$ws.Range("C4").Value = 4.660681543077796
$ws.Range("D4").Value = 4.647210928850026
$ws.Range("F4").Value = 22.1107042421702
$ws.Range("G4").Value = 3.623936959204655
$ws.Range("I4").Value = 19.11664277164639
$ws.Range("K4").Value = 7.997333942715866
$ws.Range("N4").Value = 17.56445250984624
$ws.Range("O4").Value = 19.99775884847806
$ws.Range("B5").Value = 7.994052325033453
$ws.Range("C5").Value = 4.61183598105169
$ws.Range("D5").Value = 4.622837466695643
$ws.Range("F5").Value = 22.12930233934987
$ws.Range("G5").Value = 3.624313317939431
$ws.Range("I5").Value = 19.14175583493242
$ws.Range("K5").Value = 7.94152298142639
$ws.Range("N5").Value = 17.57710673029645
$ws.Range("O5").Value = 20.02124665410441
$ws.Range("B6").Value = 7.980078490061329
$ws.Range("C6").Value = 4.603659991385159
$ws.Range("D6").Value = 4.618766462164822
$ws.Range("F6").Value = 22.13247923518634
$ws.Range("G6").Value = 3.624376506737701
$ws.Range("I6").Value = 19.14597834575961
$ws.Range("K6").Value = 7.932230732674742
$ws.Range("N6").Value = 17.57923207175832
$ws.Range("O6").Value = 20.02520910464161
$ws.Range("B7").Value = 8.076640407323001
$ws.Range("C7").Value = 4.660027198767722
$ws.Range("D7").Value = 4.646883827662274
$ws.Range("F7").Value = 22.11094911400781
$ws.Range("G7").Value = 3.623941988359157
$ws.Range("I7").Value = 19.11697793567485
$ws.Range("K7").Value = 7.996582977824502
$ws.Range("N7").Value = 17.56462155283401
$ws.Range("O7").Value = 19.99807143478623
$ws.Range("B8").Value = 8.491162910808328
$ws.Range("C8").Value = 4.898956679423945
$ws.Range("D8").Value = 4.767374752601675
$ws.Range("F8").Value = 22.02886050782544
$ws.Range("G8").Value = 3.622125301335724
$ws.Range("I8").Value = 18.99663342587872
$ws.Range("K8").Value = 8.276374118661009
$ws.Range("N8").Value = 17.50364809893972
$ws.Range("O8").Value = 19.88738508059517
$ws.Range("B9").Value = 9.254904047493685
$ws.Range("C9").Value = 5.329576413488509
$ws.Range("D9").Value = 4.989667805524442
$ws.Range("F9").Value = 21.91454767120162
$ws.Range("G9").Value = 3.618922614045816
$ws.Range("I9").Value = 18.78801553455289
$ws.Range("K9").Value = 8.805015808206395
$ws.Range("N9").Value = 17.39659837561458
$ws.Range("O9").Value = 19.70299076357522
$ws.Range("B10").Value = 9.780498552812082
$ws.Range("C10").Value = 5.620899339736531
$ws.Range("D10").Value = 5.143559366420645
$ws.Range("F10").Value = 21.85918876419506
$ws.Range("G10").Value = 3.616786805216931
$ws.Range("I10").Value = 18.65139691374564
$ws.Range("K10").Value = 9.177179797778351
$ws.Range("N10").Value = 17.32552500310853
$ws.Range("O10").Value = 19.58749599411856
$ws.Range("B11").Value = 10.01093767771096
$ws.Range("C11").Value = 5.747666028751659
$ws.Range("D11").Value = 5.211350394267756
$ws.Range("F11").Value = 21.8402502235664
$ws.Range("G11").Value = 3.615861876220393
$ws.Range("I11").Value = 18.59285445524312
$ws.Range("K11").Value = 9.342214582240878
$ws.Range("N11").Value = 17.2948252290883
$ws.Range("O11").Value = 19.53930763828057
$ws.Range("B12").Value = 10.09688809519779
$ws.Range("C12").Value = 5.794819181021203
$ws.Range("D12").Value = 5.236689488174846
$ws.Range("F12").Value = 21.83397828528125
$ws.Range("G12").Value = 3.615518304959241
$ws.Range("I12").Value = 18.57120412341249
$ws.Range("K12").Value = 9.404040266861244
$ws.Range("N12").Value = 17.28343382943363
$ws.Range("O12").Value = 19.52168689702416
$ws.Range("B13").Value = 10.07843645497301
$ws.Range("C13").Value = 5.784702065736352
$ws.Range("D13").Value = 5.23124722507843
$ws.Range("F13").Value = 21.83528902394927
$ws.Range("G13").Value = 3.615592002598266
$ws.Range("I13").Value = 18.57584385398412
$ws.Range("K13").Value = 9.390755612774891
$ws.Range("N13").Value = 17.28587677801265
$ws.Range("O13").Value = 19.52545392267875
$ws.Range("B14").Value = 10.01803551491646
$ws.Range("C14").Value = 5.751562513639487
$ws.Range("D14").Value = 5.213441781746237
$ws.Range("F14").Value = 21.83971618848064
$ws.Range("G14").Value = 3.615833476717414
$ws.Range("I14").Value = 18.59106287988583
$ws.Range("K14").Value = 9.34731473839577
$ws.Range("N14").Value = 17.29388336738341
$ws.Range("O14").Value = 19.53784539181168
$ws.Range("B15").Value = 9.980865511250117
$ws.Range("C15").Value = 5.731152184066099
$ws.Range("D15").Value = 5.202491829918789
$ws.Range("F15").Value = 21.8425451616195
$ws.Range("G15").Value = 3.615982255524737
$ws.Range("I15").Value = 18.600452484416
$ws.Range("K15").Value = 9.320617190376245
$ws.Range("N15").Value = 17.29881807911416
$ws.Range("O15").Value = 19.54551724156381
$ws.Range("B16").Value = 9.76525868245735
$ws.Range("C16").Value = 5.612496964700087
$ws.Range("D16").Value = 5.139083325963755
$ws.Range("F16").Value = 21.86055225645915
$ws.Range("G16").Value = 3.616848187826899
$ws.Range("I16").Value = 18.65529535277388
$ws.Range("K16").Value = 9.166303656820388
$ws.Range("N16").Value = 17.32756408075995
$ws.Range("O16").Value = 19.59073289433158
$ws.Range("B17").Value = 9.630723114713774
$ws.Range("C17").Value = 5.538213959874367
$ws.Range("D17").Value = 5.099607120507317
$ws.Range("F17").Value = 21.87319977686678
$ws.Range("G17").Value = 3.617391338470323
$ws.Range("I17").Value = 18.68986317170712
$ws.Range("K17").Value = 9.070503490352488
$ws.Range("N17").Value = 17.3456162624002
$ws.Range("O17").Value = 19.61958664762456
$ws.Range("B18").Value = 9.552530431793986
$ws.Range("C18").Value = 5.494947544750133
$ws.Range("D18").Value = 5.076693702460027
$ws.Range("F18").Value = 21.88106198626925
$ws.Range("G18").Value = 3.617708138131698
$ws.Range("I18").Value = 18.71008508844869
$ws.Range("K18").Value = 9.015003409073733
$ws.Range("N18").Value = 17.35615304032187
$ws.Range("O18").Value = 19.63659199102229
$ws.Range("B19").Value = 9.525918570376353
$ws.Range("C19").Value = 5.480206144262803
$ws.Range("D19").Value = 5.068900355114899
$ws.Range("F19").Value = 21.88382487099116
$ws.Range("G19").Value = 3.617816156623991
$ws.Range("I19").Value = 18.71699018226488
$ws.Range("K19").Value = 8.996145411050566
$ws.Range("N19").Value = 17.35974702239819
$ws.Range("O19").Value = 19.64241998177606
$ws.Range("B20").Value = 9.645129145164933
$ws.Range("C20").Value = 5.546177657761868
$ws.Range("D20").Value = 5.103831031128458
$ws.Range("F20").Value = 21.8717925889384
$ws.Range("G20").Value = 3.617333064639625
$ws.Range("I20").Value = 18.68614824049677
$ws.Range("K20").Value = 9.080743233220598
$ws.Range("N20").Value = 17.34367868089953
$ws.Range("O20").Value = 19.61647273118694
$ws.Range("B21").Value = 10.03581282967843
$ws.Range("C21").Value = 5.761319654471214
$ws.Range("D21").Value = 5.21868078157361
$ws.Range("F21").Value = 21.83839139472966
$ws.Range("G21").Value = 3.61576236888343
$ws.Range("I21").Value = 18.58657861383034
$ws.Range("K21").Value = 9.360092966587349
$ws.Range("N21").Value = 17.29152529465084
$ws.Range("O21").Value = 19.53418868450731
$ws.Range("B22").Value = 10.28347555636428
$ws.Range("C22").Value = 5.896960879953669
$ws.Range("D22").Value = 5.291801642950249
$ws.Range("F22").Value = 21.82180622808452
$ws.Range("G22").Value = 3.614774746362281
$ws.Range("I22").Value = 18.52452587094715
$ws.Range("K22").Value = 9.538739794458786
$ws.Range("N22").Value = 17.25880329044099
$ws.Range("O22").Value = 19.48406746742148
$ws.Range("B23").Value = 10.15201515720398
$ws.Range("C23").Value = 5.82502772035678
$ws.Range("D23").Value = 5.252957357036721
$ws.Range("F23").Value = 21.83017773675957
$ws.Range("G23").Value = 3.615298308113172
$ws.Range("I23").Value = 18.5573681233981
$ws.Range("K23").Value = 9.443768972980852
$ws.Range("N23").Value = 17.27614313691815
$ws.Range("O23").Value = 19.51048306967669
$ws.Range("B24").Value = 9.63861880993344
$ws.Range("C24").Value = 5.54257901287752
$ws.Range("D24").Value = 5.101922079204007
$ws.Range("F24").Value = 21.87242693789432
$ws.Range("G24").Value = 3.617359396109769
$ws.Range("I24").Value = 18.68782667559754
$ws.Range("K24").Value = 9.076115159596556
$ws.Range("N24").Value = 17.3445541682286
$ws.Range("O24").Value = 19.6178792342015
$ws.Range("B25").Value = 9.054172250007905
$ws.Range("C25").Value = 5.21735912696828
$ws.Range("D25").Value = 4.9311257905139
$ws.Range("F25").Value = 21.9404564785907
$ws.Range("G25").Value = 3.619750727441563
$ws.Range("I25").Value = 18.84152525357614
$ws.Range("K25").Value = 8.664583093538077
$ws.Range("N25").Value = 17.42422374510322
$ws.Range("O25").Value = 19.74937062418162
